$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows become "0M" (previously 99.99, 0.04, 316)
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# 2) Insert ten new single-value rows right after the (now) third row,
#    i.e. immediately before the row that currently holds "0".
$anchor = $t.Rows.Item(4)
$newValues = @("914", "0.00002", "0.00008", "0.00004", "0.00001", "0.00003", "0.00004", "0.00004", "0.03625", "100.0")
foreach ($val in $newValues) {
    $newRow = $t.Rows.Add($anchor)
    $newRow.Cells.Item(1).Range.Text = $val
    $anchor = $t.Rows.Item($newRow.Index + 1)
}

# 3) The last three rows (each originally containing several tab-separated
#    values crammed into one run) collapse down to a single value, reusing
#    the numbers that used to sit in rows 1-3.
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "99.99"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0.04"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "316"
